$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.592.32'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '2.246.17'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.45'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.60'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.86'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0951'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.23'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.55'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.861'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '2.260.44'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '42.319.63'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000103'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.86%  '
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.04'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.84'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +49.36%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.77'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.73'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.64%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.31'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +4.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.17'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.72'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.09'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +23.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0820'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.46'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -9.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.126'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.67'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0316'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '13.57'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.91%  '
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.72'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '63.69'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.97%  '
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.14'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.87'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.78%  '
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.43'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.15'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.18'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("E51").Value = '  +1.38%  '
